# Episode Choice Remake 3.0
# Add 5 new games to the "game time" tracking sheet, then re-sort the whole
# data range (A2:D) descending by column D ("Время, ч" / hours played), and
# reset the sheet view back to the top with D15 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Append the 5 new rows of data right after the current last row (64) ---
# Column A = game title, column B = episode count (left blank, same as many
# existing rows), column C = total episode/run count, column D = hours.

$newGames = @(
    @("Atomic Heart", 30, 22.4166666666666),
    @("S.T.A.L.K.E.R.: Clear Sky", 16, 12.5),
    @("S.T.A.L.K.E.R.: Lost Alpha Enhanced Edition [ng++]", 4, 7.46666666666666),
    @("Dead Space", 13, 11.3166666666666),
    @("S.T.A.L.K.E.R.: Call of Pripyat Gunslinger", 17, 16)
)

$startRow = 64
for ($i = 0; $i -lt $newGames.Count; $i++) {
    $r = $startRow + $i
    $game = $newGames[$i]
    $ws.Cells.Item($r, 1).Value = $game[0]
    $ws.Cells.Item($r, 3).Value = $game[1]
    $ws.Cells.Item($r, 4).Value = $game[2]
}

# --- 2. Re-sort the whole data range (A2:D68) descending by column D ---
$lastRow = 68
$sortRange = $ws.Range("A2:D$lastRow")
$keyRange = $ws.Range("D2:D$lastRow")
$sortRange.Sort($keyRange, 2)

# --- 3. Reset the sheet view: scroll back to the top, select D15 ---
$ws.Range("D15").Select()
